# Auto-generated Excel COM-interop script to apply scheduled market-data refresh
# to the Goblin Profits workbook (Leve crafting profit tracker).
# For every affected row, currentAveragePrice* / LevePrice* / LeveProfit* columns
# (H..N) are updated to the latest computed snapshot values. A couple of rows also
# gain or lose a trailing LeveProfitHQ (N) or LeveProfitNQ (M) cell because the
# refreshed computation changed whether an HQ/NQ variant applies to that leve.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2218.6365
$ws.Range("I137").Value = 2360.6
$ws.Range("J137").Value = 2100.3333
$ws.Range("K137").Value = 7081.799999999999
$ws.Range("L137").Value = 6300.999899999999
$ws.Range("M137").Value = -4531.799999999999
$ws.Range("N137").Value = -11400.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 499.5
$ws.Range("I2").Value = 499.5
$ws.Range("K2").Value = 499.5
$ws.Range("M2").Value = -386.5

$ws.Range("H23").Value = 22503.25
$ws.Range("J23").Value = 22503.25
$ws.Range("L23").Value = 22503.25
$ws.Range("N23").Value = -23021.25

$ws.Range("H26").Value = 1227
$ws.Range("I26").Value = 1227
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 1227
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -897
$ws.Range("N26").ClearContents()

$ws.Range("H63").Value = 8599.799999999999
$ws.Range("J63").Value = 9999.666999999999
$ws.Range("L63").Value = 9999.666999999999
$ws.Range("N63").Value = -11371.667

$ws.Range("H66").Value = 8599.799999999999
$ws.Range("J66").Value = 9999.666999999999
$ws.Range("L66").Value = 49998.335
$ws.Range("N66").Value = -56862.335

$ws.Range("H102").Value = 4885.643
$ws.Range("I102").Value = 2779.1052
$ws.Range("J102").Value = 9332.777
$ws.Range("K102").Value = 2779.1052
$ws.Range("L102").Value = 9332.777
$ws.Range("M102").Value = -1157.1052
$ws.Range("N102").Value = -12576.777

$ws.Range("H109").Value = 55000
$ws.Range("J109").Value = 55000
$ws.Range("L109").Value = 55000
$ws.Range("N109").Value = -57774

$ws.Range("H116").Value = 499.5
$ws.Range("I116").Value = 499.5
$ws.Range("K116").Value = 499.5
$ws.Range("M116").Value = 1794.5

$ws.Range("H132").Value = 1586.2142
$ws.Range("I132").Value = 1759.7727
$ws.Range("J132").Value = 949.8333
$ws.Range("K132").Value = 5279.3181
$ws.Range("L132").Value = 2849.4999
$ws.Range("M132").Value = -2749.3181
$ws.Range("N132").Value = -7909.4999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 499.5
$ws.Range("I3").Value = 499.5
$ws.Range("K3").Value = 499.5
$ws.Range("M3").Value = -385.5

$ws.Range("H20").Value = 1779.7333
$ws.Range("I20").Value = 1299.8
$ws.Range("K20").Value = 1299.8
$ws.Range("M20").Value = -1052.8

$ws.Range("H94").Value = 1066.3158
$ws.Range("I94").Value = 1106.1177
$ws.Range("K94").Value = 1106.1177
$ws.Range("M94").Value = -655.1177

$ws.Range("H99").Value = 2830.5652
$ws.Range("I99").Value = 1390.8
$ws.Range("J99").Value = 3938.077
$ws.Range("K99").Value = 1390.8
$ws.Range("L99").Value = 3938.077
$ws.Range("M99").Value = 107.2
$ws.Range("N99").Value = -6934.077

$ws.Range("H134").Value = 2808.261
$ws.Range("I134").Value = 2736.875
$ws.Range("K134").Value = 8210.625
$ws.Range("M134").Value = -5675.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1928.5714
$ws.Range("I16").Value = 1333.3334
$ws.Range("J16").Value = 2375
$ws.Range("K16").Value = 1333.3334
$ws.Range("L16").Value = 2375
$ws.Range("M16").Value = -1046.3334
$ws.Range("N16").Value = -2949

$ws.Range("H58").Value = 2161.5217
$ws.Range("I58").Value = 2145.2942
$ws.Range("K58").Value = 2145.2942
$ws.Range("M58").Value = -1942.2942

$ws.Range("H94").Value = 3118.0625
$ws.Range("I94").Value = 2429.2856
$ws.Range("J94").Value = 3653.7778
$ws.Range("K94").Value = 2429.2856
$ws.Range("L94").Value = 3653.7778
$ws.Range("M94").Value = -1978.2856
$ws.Range("N94").Value = -4555.7778

$ws.Range("H99").Value = 2968.5
$ws.Range("J99").Value = 2500
$ws.Range("L99").Value = 2500
$ws.Range("N99").Value = -5496

$ws.Range("H113").Value = 1928.5714
$ws.Range("I113").Value = 1333.3334
$ws.Range("J113").Value = 2375
$ws.Range("K113").Value = 1333.3334
$ws.Range("L113").Value = 2375
$ws.Range("M113").Value = 836.6666
$ws.Range("N113").Value = -6715

$ws.Range("H126").Value = 2968.5
$ws.Range("J126").Value = 2500
$ws.Range("L126").Value = 7500
$ws.Range("N126").Value = -12440

$ws.Range("H136").Value = 2161.5217
$ws.Range("I136").Value = 2145.2942
$ws.Range("K136").Value = 6435.882599999999
$ws.Range("M136").Value = -3885.882599999999

$ws.Range("H141").Value = 161249.75
$ws.Range("I141").Value = 15000
$ws.Range("K141").Value = 15000
$ws.Range("M141").Value = -9820

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1319.9
$ws.Range("I23").Value = 99.5
$ws.Range("J23").Value = 1625
$ws.Range("K23").Value = 298.5
$ws.Range("L23").Value = 4875
$ws.Range("M23").Value = -63.5
$ws.Range("N23").Value = -5345

$ws.Range("H35").Value = 200
$ws.Range("I35").Value = 200
$ws.Range("K35").Value = 600
$ws.Range("M35").Value = -312

$ws.Range("H112").Value = 495
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H127").Value = 2148.5
$ws.Range("J127").Value = 2148.5
$ws.Range("L127").Value = 6445.5
$ws.Range("N127").Value = -16365.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 750
$ws.Range("I5").Value = 750
$ws.Range("K5").Value = 750
$ws.Range("M5").Value = -638

$ws.Range("H70").Value = 9378.571
$ws.Range("J70").Value = 9461.166999999999
$ws.Range("L70").Value = 9461.166999999999
$ws.Range("N70").Value = -10001.167

$ws.Range("H73").Value = 9378.571
$ws.Range("J73").Value = 9461.166999999999
$ws.Range("L73").Value = 9461.166999999999
$ws.Range("N73").Value = -11333.167

$ws.Range("H80").Value = 3745
$ws.Range("I80").Value = 2022.4166
$ws.Range("K80").Value = 2022.4166
$ws.Range("M80").Value = -1024.4166

$ws.Range("H83").Value = 3745
$ws.Range("I83").Value = 2022.4166
$ws.Range("K83").Value = 10112.083
$ws.Range("M83").Value = -5120.083000000001

$ws.Range("H132").Value = 3391.9333
$ws.Range("I132").Value = 2988.182
$ws.Range("J132").Value = 4502.25
$ws.Range("K132").Value = 8964.545999999998
$ws.Range("L132").Value = 13506.75
$ws.Range("M132").Value = -6434.545999999998
$ws.Range("N132").Value = -18566.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 23069.2
$ws.Range("I53").Value = 20673
$ws.Range("J53").Value = 24666.666
$ws.Range("K53").Value = 20673
$ws.Range("L53").Value = 24666.666
$ws.Range("M53").Value = -20155
$ws.Range("N53").Value = -25702.666

$ws.Range("H122").Value = 7960.478
$ws.Range("I122").Value = 6122.8667
$ws.Range("K122").Value = 18368.6001
$ws.Range("M122").Value = -15918.6001

$ws.Range("H130").Value = 70000
$ws.Range("J130").Value = 70000
$ws.Range("L130").Value = 70000
$ws.Range("N130").Value = -80040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3203.5
$ws.Range("I81").Value = 938
$ws.Range("K81").Value = 1876
$ws.Range("M81").Value = -815

$ws.Range("H84").Value = 3203.5
$ws.Range("I84").Value = 938
$ws.Range("K84").Value = 9380
$ws.Range("M84").Value = -4076

$ws.Range("H93").Value = 75000
$ws.Range("I93").Value = 62500
$ws.Range("J93").Value = 100000
$ws.Range("K93").Value = 62500
$ws.Range("L93").Value = 100000
$ws.Range("M93").Value = -60004
$ws.Range("N93").Value = -104992

$ws.Range("H96").Value = 4111
$ws.Range("J96").Value = 4000
$ws.Range("L96").Value = 4000
$ws.Range("N96").Value = -6746

$ws.Range("H107").Value = 1949.6666
$ws.Range("I107").Value = 1834.3077
$ws.Range("J107").Value = 2699.5
$ws.Range("K107").Value = 5502.9231
$ws.Range("L107").Value = 8098.5
$ws.Range("M107").Value = -3582.9231
$ws.Range("N107").Value = -11938.5

$ws.Range("H113").Value = 2488.2
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 6000
$ws.Range("M113").Value = -3830

$ws.Range("H114").Value = 78996
$ws.Range("J114").Value = 78996
$ws.Range("L114").Value = 78996
$ws.Range("N114").Value = -87674

$ws.Range("H132").Value = 2980.147
$ws.Range("I132").Value = 2891.2258
$ws.Range("K132").Value = 8673.6774
$ws.Range("M132").Value = -6143.6774

$ws.Range("H136").Value = 1636.5143
$ws.Range("I136").Value = 1116.4615
$ws.Range("K136").Value = 3349.3845
$ws.Range("M136").Value = -799.3844999999997
